$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "2020-08-28T00:00:00.000Z"
$ws.Range("B4").Value = "ZIZ21"
$ws.Range("C4").Value = "ICE Silver 5000oz Dec21"
$ws.Range("E4").Value = "Commodities"
$ws.Range("F4").Value = 24.3
$ws.Range("G4").Value = 25.4
$ws.Range("H4").Value = "imd_12457801"

$ws.Range("K10").Select()
